$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '27.930.56'
Set-TextValue 'E2' '  +0.71%  '
Set-TextValue 'D3' '1.882.25'
Set-TextValue 'E3' '  +0.32%  '
Set-TextValue 'D4' '1.018'
Set-TextValue 'E4' '  +1.44%  '
Set-TextValue 'D5' '335.16'
Set-TextValue 'E5' '  +0.72%  '
Set-TextValue 'E6' '  +1.37%  '
Set-TextValue 'D7' '0.4681'
Set-TextValue 'E7' '  -1.10%  '
Set-TextValue 'D8' '0.3913'
Set-TextValue 'E8' '  -1.21%  '
Set-TextValue 'D9' '47.03'
Set-TextValue 'E9' '  -1.51%  '
Set-TextValue 'D10' '0.07949'
Set-TextValue 'E10' '  -1.17%  '
Set-TextValue 'D11' '1.009'
Set-TextValue 'E11' '  -1.55%  '
Set-TextValue 'D12' '21.63'
Set-TextValue 'E12' '  -1.23%  '
Set-TextValue 'D13' '1.906.58'
Set-TextValue 'E13' '  +1.67%  '
Set-TextValue 'D14' '5.947'
Set-TextValue 'E14' '  -0.23%  '
Set-TextValue 'D15' '7.100'
Set-TextValue 'E15' '  -0.91%  '
Set-TextValue 'D16' '1.019'
Set-TextValue 'E16' '  +1.30%  '
Set-TextValue 'D17' '0.06794'
Set-TextValue 'E17' '  +2.40%  '
Set-TextValue 'E18' '  +0.20%  '
Set-TextValue 'D19' '0.00001044'
Set-TextValue 'E19' '  -0.70%  '
Set-TextValue 'D20' '16.98'
Set-TextValue 'E20' '  -1.59%  '
Set-TextValue 'E21' '  +1.44%  '
Set-TextValue 'D22' '27.936.09'
Set-TextValue 'E22' '  +0.73%  '
Set-TextValue 'D23' '5.459'
Set-TextValue 'E23' '  -0.69%  '
Set-TextValue 'D24' '10.96'
Set-TextValue 'E24' '  -0.86%  '
Set-TextValue 'E25' '  +2.56%  '
Set-TextValue 'D26' '2.133.95'
Set-TextValue 'E26' '  +1.71%  '
Set-TextValue 'D27' '159.22'
Set-TextValue 'E27' '  +1.63%  '
Set-TextValue 'D28' '19.98'
Set-TextValue 'E28' '  -1.20%  '
Set-TextValue 'D29' '2.065'
Set-TextValue 'E29' '  -1.85%  '
Set-TextValue 'D30' '5.446'
Set-TextValue 'E30' '  -2.46%  '
Set-TextValue 'D31' '120.70'
Set-TextValue 'E31' '  -1.54%  '
Set-TextValue 'D32' '0.09520'
Set-TextValue 'E32' '  -0.44%  '
Set-TextValue 'D33' '0.9555'
Set-TextValue 'E33' '  -1.58%  '
Set-TextValue 'E34' '  +0.75%  '
Set-TextValue 'D35' '5.312'
Set-TextValue 'E35' '  +0.11%  '
Set-TextValue 'D36' '1.348'
Set-TextValue 'E36' '  -7.41%  '
Set-TextValue 'D37' '0.06114'
Set-TextValue 'E37' '  +0.03%  '
Set-TextValue 'D38' '0.02234'
Set-TextValue 'E38' '  -1.45%  '
Set-TextValue 'D39' '1.205'
Set-TextValue 'E39' '  -1.83%  '
Set-TextValue 'E40' '  +1.39%  '
Set-TextValue 'D41' '8.113'
Set-TextValue 'E41' '  -0.90%  '
Set-TextValue 'D42' '0.5863'
Set-TextValue 'E42' '  -2.26%  '
Set-TextValue 'E43' '  -1.15%  '
Set-TextValue 'E44' '  -1.10%  '
Set-TextValue 'E45' '  +0.40%  '
Set-TextValue 'E46' '  -1.70%  '
Set-TextValue 'D47' '12.14'
Set-TextValue 'E47' '  -1.61%  '
Set-TextValue 'D48' '3.397'
Set-TextValue 'E48' '  -0.48%  '
Set-TextValue 'E49' '  -1.13%  '
Set-TextValue 'E50' '  +0.64%  '
Set-TextValue 'D51' '113.31'
Set-TextValue 'E51' '  +0.42%  '
